{"js": "// Remove the \"line: <N>\" paragraph from the test-log table (it immediately\n// follows the \"header: ...\" paragraph in the same table cell). Office test\n// logs either drop the line-number paragraph entirely (this case) or gain\n// one; here it is being dropped.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text === \"line: 4\") {\n    toDelete.push(paragraphs.items[i]);\n  }\n}\n\nfor (const paragraph of toDelete) {\n  paragraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"line: <N>\" paragraph from the test-log table (it immediately\n# follows the \"header: ...\" paragraph in the same table cell). Office test\n# logs either drop the line-number paragraph entirely (this case) or gain\n# one; here it is being dropped.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    # A paragraph that is the last one in a table cell ends its Range.Text\n    # with CR (13) + BEL (7) instead of just CR (13); strip both before\n    # comparing so the match works regardless of position in the cell.\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"line: 4\") {\n        $p.Range.Delete()\n    }\n}\n"}
